# Updates the cryptos list (Price + Volume(1h) columns) for rows 2-51
# per the "Updated cryptos list on Mon May  8 11:38:56 UTC 2023 with GitHub Actions" commit.
#
# Column D (Price) holds plain text such as "27.979.14" or "1.001" - these are
# *not* real numbers (some even contain two "." separators), they must stay as
# text. Excel auto-converts a typed value like "1.001" into the number 1.001,
# so for any new Price value that COM would otherwise misinterpret as numeric we
# force the cell to Text format first, then assign the literal string. Values that
# are safely non-numeric (two dots, e.g. "27.979.14") are assigned directly.
#
# Column E (Volume(1h)) is always a padded percentage string (e.g. "  -3.12%  ")
# and is never numeric, so it is always assigned directly.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '27.979.14'
$ws.Range("E2").Value = '  -3.12%  '

$ws.Range("D3").Value = '1.864.68'
$ws.Range("E3").Value = '  -2.14%  '

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.001'
$ws.Range("E4").Value = '  -0.15%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '318.32'
$ws.Range("E5").Value = '  -1.93%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '1.000'
$ws.Range("E6").Value = '  -0.09%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.4374'
$ws.Range("E7").Value = '  -4.60%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.3704'
$ws.Range("E8").Value = '  -2.87%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.07516'
$ws.Range("E9").Value = '  -2.58%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.9384'
$ws.Range("E10").Value = '  -4.11%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '21.32'
$ws.Range("E11").Value = '  -3.97%  '

$ws.Range("D12").Value = '1.886.57'
$ws.Range("E12").Value = '  +0.35%  '

$ws.Range("E13").Value = '  -3.00%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '5.452'
$ws.Range("E14").Value = '  -3.91%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.06824'
$ws.Range("E15").Value = '  -3.42%  '

$ws.Range("E16").Value = '  -0.15%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '81.72'
$ws.Range("E17").Value = '  -2.54%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.000009050'
$ws.Range("E18").Value = '  -4.24%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '1.000'
$ws.Range("E19").Value = '  -0.05%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '15.95'
$ws.Range("E20").Value = '  -4.11%  '

$ws.Range("D21").Value = '27.947.82'
$ws.Range("E21").Value = '  -3.19%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '5.113'
$ws.Range("E22").Value = '  -3.71%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '11.08'
$ws.Range("E23").Value = '  +1.35%  '

$ws.Range("D24").Value = '2.094.12'
$ws.Range("E24").Value = '  -0.64%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.004'
$ws.Range("E25").Value = '  -4.59%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '154.34'
$ws.Range("E26").Value = '  -2.82%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '18.41'
$ws.Range("E27").Value = '  -3.12%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '5.435'
$ws.Range("E28").Value = '  -4.21%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '113.70'
$ws.Range("E29").Value = '  -3.28%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '1.735'
$ws.Range("E30").Value = '  -7.38%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '0.08998'
$ws.Range("E31").Value = '  -3.28%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.8128'
$ws.Range("E32").Value = '  -5.67%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '4.826'
$ws.Range("E33").Value = '  -5.32%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '1.177'
$ws.Range("E34").Value = '  -5.33%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '2.936'
$ws.Range("E35").Value = '  -3.07%  '

$ws.Range("E36").Value = '  -0.01%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.05497'
$ws.Range("E37").Value = '  -3.69%  '

$ws.Range("E38").Value = '  -3.68%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.01981'
$ws.Range("E39").Value = '  -2.95%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '2.897'
$ws.Range("E40").Value = '  -0.48%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.5271'
$ws.Range("E41").Value = '  -3.86%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '7.057'
$ws.Range("E42").Value = '  -5.35%  '

$ws.Range("E43").Value = '  -3.50%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '8.817'
$ws.Range("E44").Value = '  -5.53%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.06784'
$ws.Range("E45").Value = '  -1.46%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.4909'
$ws.Range("E46").Value = '  -5.12%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '10.66'
$ws.Range("E47").Value = '  -5.33%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '106.50'
$ws.Range("E48").Value = '  -3.57%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '1.683'
$ws.Range("E49").Value = '  -5.29%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '1.000'
$ws.Range("E50").Value = '  -0.14%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '1.904'
$ws.Range("E51").Value = '  -12.23%  '
